$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The results table is built from heavily vMerge'd cells, and several of the
# target cells span multiple paragraphs (two blank spacer paragraphs followed
# by the real text paragraph). Because of that we address every edit by the
# absolute paragraph index within the table's Range.Paragraphs collection
# (which reports correct, document-order Start/End offsets) rather than by
# Table.Cell()/Find, which are unreliable across multi-paragraph cells and
# ambiguous for short repeated numbers like "89".
$paras = $t.Range.Paragraphs

# index (1-based, within $t.Range.Paragraphs) -> new text
$edits = @{
    13  = "Sun Microsystem Pvt Ltd";   # "Sun Microsystem Pvt Ltd - 90-9-"
    14  = "77908954";                  # "77908954090"
    16  = "Chennai";                   # "Chennai-98"
    19  = "Cuddalore";                 # "Cuddalore-1"
    21  = "18-11-2023 12.00.00 AM";    # "10-11-2023 12.00.00 AM"
    29  = "adada";                     # "Ammenometer"
    30  = "adad";                      # "12345678"
    31  = "18-11-2023 12.00.00 AM";    # "03-11-2023 12.00.00 AM"
    34  = "adad adad";                 # "12345678 12345678"
    35  = "07-11-2023 12.00.00 AM";    # "16-11-2023 12.00.00 AM"
    76  = "RSVD-2";                    # "RSVD-3"            (row 10, S3)
    81  = "0.25";                      # "0.998"
    82  = "1896";                      # "1898"
    83  = "1856";                      # "198"
    84  = "4856";                      # "8198"
    85  = "5698";                      # "189"
    86  = "9874";                      # "81981"
    87  = "4836";                      # "18493"
    88  = "1209";                      # "18456"
    91  = "4670";                      # "3512232"
    94  = "2500";                      # "2078"
    97  = "112";                       # "101412"
    102 = "0.56";                      # "176"               (row 11, S4)
    103 = "1890";                      # "8989"
    104 = "9009";                      # "89189"
    105 = "1909";                      # "89"
    106 = "9010";                      # "89"
    107 = "9080";                      # "898"
    108 = "6180";                      # "19851"
    109 = "3461";                      # "3493776"
}

# Snapshot the current Start/End for every paragraph we're about to touch
# *before* making any edits (editing shifts every later offset).
$targets = @()
foreach ($idx in $edits.Keys) {
    $p = $paras.Item($idx)
    $targets += [PSCustomObject]@{
        Start = $p.Range.Start
        End   = $p.Range.End
        New   = $edits[$idx]
    }
}

# Apply from the end of the document backwards so earlier (still queued)
# offsets stay valid.
$targets = $targets | Sort-Object -Property Start -Descending
foreach ($item in $targets) {
    $rng = $d.Range($item.Start, $item.End - 1)
    $rng.Text = $item.New
}

# Finally, drop the whole duplicated RSVD-3 / NLT-20 row (the "S3" row that
# repeats after the "S4" row) now that every surviving cell has its new text.
$t.Rows.Item(12).Delete()
